# iti bug fixed. itis changed
# The inter-trial-interval (iti) values in column A (rows 2-41) were
# recalculated. The previous repeating cycle of 4.5 / 6 / 7.5 / 9 seconds
# is replaced by a corrected repeating cycle of 3 / 5 / 7 / 9 seconds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$itiValues = @(3, 5, 7, 9)

$row = 2
for ($i = 0; $i -lt 10; $i++) {
    foreach ($val in $itiValues) {
        $ws.Cells.Item($row, 1).Value = $val
        $row++
    }
}

# Restore the scroll position of the sheet view (was showing row 24 at the
# top, now shows row 13 at the top), keeping the selection on A41.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("A41").Select() | Out-Null

# Re-localize the default cell style name to the English "Normal"
# (it was saved previously as the German "Standard").
try {
    $style = $wb.Styles.Item("Standard")
    $style.Name = "Normal"
} catch {
    # Style already named "Normal" or not available - ignore.
}
